# Split the run "a unified means of access control for data. " into two
# runs: "a unified means of access control for " and "data" (dropping the
# trailing ". "), matching the other bullet's "... analysing " + "data"
# pattern on the same slide.

$p = $ppt.ActivePresentation

$oldRunText  = "a unified means of access control for data. "
$firstPart   = "a unified means of access control for "
$secondPart  = "data"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if (-not $sh.HasTextFrame) { continue }
        if (-not $sh.TextFrame.HasText) { continue }

        $tr = $sh.TextFrame.TextRange
        $paraCount = $tr.Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            $runCount = $para.Count

            for ($ri = 1; $ri -le $runCount; $ri++) {
                $run = $para.Runs($ri, 1)

                if ($run.Text -eq $oldRunText) {
                    # Shrink the existing run to the first half of the text,
                    # then add a brand-new run right after it holding "data".
                    $run.Text = $firstPart
                    $run.InsertAfter($secondPart) | Out-Null
                }
            }
        }
    }
}
